$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B to hold "mu"
# (shifts old B..F -> C..G, matching the target layout)
$ws.Columns.Item(2).Insert()

# New header + data: mu = 2 for every row
$ws.Range("B1").Value = "mu"
$ws.Range("B2:B8").Value = 2

# One-decimal number format for the mu column (incl. header cell, which
# also ends up carrying the style in the source workbook)
$ws.Range("B1:B8").NumberFormat = "0.0"

# Drop the hidden chart-linked defined names (_xlchart.v1.*) that pointed
# at the pre-edit column layout
while ($wb.Names.Count -gt 0) {
    $wb.Names.Item(1).Delete()
}
